$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Add a new "Weiteres" header in O1, copying the formatting (style) that
# N1 ("Teilnehmeranzahl") currently has, then relabel the header cells so
# the whole header row reads:
# Jahr | Erstautor | Typ | Space | Entfernungen | Messverfahren |
# Displaytyp | VE | Titel | Untersuchungsgegenstand | Ergebnisse |
# Experiment-Details | Teilnehmeranzahl | Trials | Weiteres
$ws.Range("N1").Copy($ws.Range("O1"))

$ws.Range("B1").Value = "Erstautor"
$ws.Range("D1").Value = "Space"
$ws.Range("M1").Value = "Teilnehmeranzahl"
$ws.Range("N1").Value = "Trials"
$ws.Range("O1").Value = "Weiteres"

# --- Data rows (Piryankova literature entries) -------------------------
$ws.Range("A2").Value = 2013
$ws.Range("B2").Value = "Piryankova"
$ws.Range("C2").Value = "E"
$ws.Range("D2").Value = "A"
$ws.Range("E2").Value = "1,5-6m"
$ws.Range("M2").Value = 77
$ws.Range("N2").Value = "27-30"

$ws.Range("C3").Value = "E"
$ws.Range("D3").Value = "A"
$ws.Range("E3").Value = "2-6m"
$ws.Range("L3").Value = "RW"
$ws.Range("M3").Value = 16

$ws.Range("E4").Value = "1,5-5,5m"
$ws.Range("L4").Value = "semi-spherical LSID"
$ws.Range("M4").Value = 11

$ws.Range("E5").Value = "2-5,5m"
$ws.Range("L5").Value = "MPI cabin"
$ws.Range("M5").Value = 10

$ws.Range("E6").Value = "2-5,5m"
$ws.Range("L6").Value = "flat LSID"
$ws.Range("M6").Value = 40

# --- Selection matches the authored workbook ---------------------------
$ws.Range("E7").Select()
